# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as a new row 58 (pushing the
# existing rows 58..123 down to 59..124). The rest of the data is unchanged;
# only the table grows by one row (A1:R123 -> A1:R124).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 58, shifting rows 58-123 down to 59-124.
$ws.Rows.Item(58).Insert()

# Populate the newly inserted row 58 with the new observation.
$ws.Cells.Item(58, 1).Value = 7
$ws.Cells.Item(58, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(58, 3).Value = "Ñuble"
$ws.Cells.Item(58, 4).Value = 45195
$ws.Cells.Item(58, 5).Value = 16
$ws.Cells.Item(58, 6).Value = 100112001
$ws.Cells.Item(58, 7).Value = "Berenjena"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 50
$ws.Cells.Item(58, 11).Value = 10000
$ws.Cells.Item(58, 12).Value = 10000
$ws.Cells.Item(58, 13).Value = 10000
$ws.Cells.Item(58, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(58, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(58, 16).Value = 167
$ws.Cells.Item(58, 17).Value = 60
$ws.Cells.Item(58, 18).Value = "Hortaliza"

# Give the new date cell the same date/time number format used by the rest
# of column D (style index 2 in the original file), matching what Insert()
# already carried down from the old row 58.
$ws.Cells.Item(58, 4).NumberFormat = $ws.Cells.Item(59, 4).NumberFormat
